# "Path to Graduation" planner update.
# The sheet is a repeating 3-column (Fall/Spring/Summer) x N-year grid of
# semester course tables. This edit:
#   - reshuffles/renames several course rows within the existing
#     2022-2023 and 2023-2024 blocks (new CYBR/CPSC course codes),
#   - adds two new data rows (16 & 17) to the 2023-2024 Fall column,
#   - trims the 2022-2023 block's Fall column back down to rows 4-7
#     (row 8 keeps only the Spring entry now),
#   - replaces the lone 2024-2025 Fall row with two new rows (22 & 23),
#   - appends two brand-new year blocks: Fall/Spring/Summer 2026
#     (rows 39-47) and Fall/Spring/Summer 2027 (rows 48-56), each with
#     a header row and a Total row with SUM formulas, mirroring the
#     existing blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2022-2023 block (rows 4-11) tweaks ---
$ws.Range("C4").Value = "CYBR 3115"
$ws.Range("C5").Value = "CYBR 2159"
$ws.Range("A6").Value = "CYBR 3119"
$ws.Range("C6").Value = "CPSC 1302"
$ws.Range("A7").Value = "CYBR 3106"
$ws.Range("C7").Value = "CPSC 2108"
$ws.Range("D7").Value = 3
$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "CYBR 3108"
$ws.Range("D8").Value = 3

# --- 2023-2024 block (rows 12-20) tweaks, plus two new rows ---
$ws.Range("A13").Value = "CPSC 4155"
$ws.Range("C13").Value = "CPSC 4135"
$ws.Range("A14").Value = "DSCI 3111"
$ws.Range("C14").Value = "CPSC 4175"
$ws.Range("A15").Value = "CPSC 3165"
$ws.Range("C15").Value = "CPSC 6180"
$ws.Range("A16").Value = "CPSC 4111"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "CPSC 6185"
$ws.Range("D16").Value = 3
$ws.Range("A17").Value = "CPSC 4148"
$ws.Range("B17").Value = 3

# --- 2024-2025 block (rows 21-29) tweaks ---
$ws.Range("A22").Value = "CPSC 6985"
$ws.Range("B22").Value = 4
$ws.Range("A23").Value = "CPSC 4000"
$ws.Range("B23").Value = 0

# --- New Fall/Spring/Summer 2026 block (rows 39-47) ---
$ws.Range("A39").Value = "Fall 2026"
$ws.Range("B39").Value = "Credits"
$ws.Range("C39").Value = "Spring 2026"
$ws.Range("D39").Value = "Credits"
$ws.Range("E39").Value = "Summer 2026"
$ws.Range("F39").Value = "Credits"

$ws.Range("A47").Value = "Total"
$ws.Range("B47").Formula = "=SUM(B40:B46)"
$ws.Range("C47").Value = "Total"
$ws.Range("D47").Formula = "=SUM(D40:D46)"
$ws.Range("E47").Value = "Total"
$ws.Range("F47").Formula = "=SUM(F40:F46)"

# --- New Fall/Spring/Summer 2027 block (rows 48-56) ---
$ws.Range("A48").Value = "Fall 2027"
$ws.Range("B48").Value = "Credits"
$ws.Range("C48").Value = "Spring 2027"
$ws.Range("D48").Value = "Credits"
$ws.Range("E48").Value = "Summer 2027"
$ws.Range("F48").Value = "Credits"

$ws.Range("A56").Value = "Total"
$ws.Range("B56").Formula = "=SUM(B49:B55)"
$ws.Range("C56").Value = "Total"
$ws.Range("D56").Formula = "=SUM(D49:D55)"
$ws.Range("E56").Value = "Total"
$ws.Range("F56").Formula = "=SUM(F49:F55)"
